$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")
$ws.Columns("O:O").Insert()
